# Added support for local execution videos
$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# --- "Count" sheet: B3 flips from TRUE to FALSE; selection moves from B2 to B3 ---
$wsCount = $wb.Worksheets.Item("Count")
$wsCount.Range("B3").Value = $false

# --- "TestData" sheet: row 2 browser/mode change + new Version values + F5 flip ---
$wsData = $wb.Worksheets.Item("TestData")

$wsData.Range("C2").Value = "chrome"
$wsData.Range("G2").Value = "remote"

# Column D ("Version") goes from a blank quote-prefixed text cell to the number 109,
# while keeping the cell's existing (quote-prefix) style. Setting .Value alone drops
# that style, so re-apply the original cell format afterwards by pasting formats from
# an equivalent untouched style-4 cell in the same column family.
$wsData.Range("D2").Value = 109
$wsData.Range("J3").Copy()
$wsData.Range("D2").PasteSpecial($xlPasteFormats)

$wsData.Range("D3").Value = 109
$wsData.Range("J4").Copy()
$wsData.Range("D3").PasteSpecial($xlPasteFormats)

$wsData.Range("D4").Value = 109
$wsData.Range("H5").Copy()
$wsData.Range("D4").PasteSpecial($xlPasteFormats)

$wsData.Range("D5").Value = 109
$wsData.Range("J5").Copy()
$wsData.Range("D5").PasteSpecial($xlPasteFormats)

$wsData.Range("F5").Value = $false

$excel.CutCopyMode = $false

# --- Update the selections shown in each sheet's view, restoring the active tab after ---
$wsCount.Activate()
$wsCount.Range("B3").Select()

$wsData.Activate()
$wsData.Range("D2:D5").Select()
